$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "qwen2:7b-instruct-q5_K_M"
$ws.Range("B26").Value = "llama3:70b"
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 200
$ws.Range("E26").Value = 2531.55
$ws.Range("F26").Value = 361.15
$ws.Range("G26").Value = 0.3375
$ws.Range("H26").Value = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_5_200_test_match.txt"
$ws.Range("I26").Value = 361.15
$ws.Range("J26").Value = 0.9
$ws.Range("K26").Value = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_5_200_test_correct.txt"
$ws.Range("L26").Value = 361.15
$ws.Range("M26").Value = 0.7375
$ws.Range("N26").Value = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_5_200_test_executable.txt"
$ws.Range("O26").Value = 0
$ws.Range("P26").Value = 0
$ws.Range("Q26").Value = 362.48
$ws.Range("R26").Value = 0.3625
$ws.Range("S26").Value = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_5_200_test_fewshot_match.txt"
$ws.Range("T26").Value = 362.48
$ws.Range("U26").Value = 0.825
$ws.Range("V26").Value = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_5_200_test_fewshot_correct.txt"
$ws.Range("W26").Value = 362.48
$ws.Range("X26").Value = 0.6875
$ws.Range("Y26").Value = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_5_200_test_fewshot_executable.txt"
$ws.Range("Z26").Value = 0
$ws.Range("AA26").Value = 1370.61
$ws.Range("AB26").Value = 437.32
$ws.Range("AC26").Value = 0.1125
$ws.Range("AD26").Value = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_5_200_test_bootstrap_match.txt"
$ws.Range("AE26").Value = 437.32
$ws.Range("AF26").Value = 0.7125
$ws.Range("AG26").Value = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_5_200_test_bootstrap_correct.txt"
$ws.Range("AH26").Value = 437.32
$ws.Range("AI26").Value = 0.6375
$ws.Range("AJ26").Value = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_5_200_test_bootstrap_executable.txt"
$ws.Range("AK26").Value = 0
$ws.Range("AL26").Value = 2
$ws.Range("AM26").Value = 2
